$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-08-23 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-24 Saturday", 2) | Out-Null

$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "29+52="
$t.Cell(1, 2).Range.Text = "12-5="
$t.Cell(1, 3).Range.Text = "54+40="
$t.Cell(1, 4).Range.Text = "73-39="
$t.Cell(1, 5).Range.Text = "45-44="

$t.Cell(2, 1).Range.Text = "47+26="
$t.Cell(2, 2).Range.Text = "37+32="
$t.Cell(2, 3).Range.Text = "6-0="
$t.Cell(2, 4).Range.Text = "57+35="
$t.Cell(2, 5).Range.Text = "72+19="

$t.Cell(3, 1).Range.Text = "88-0="
$t.Cell(3, 2).Range.Text = "11+60="
$t.Cell(3, 3).Range.Text = "16+75="
$t.Cell(3, 4).Range.Text = "38-12="
$t.Cell(3, 5).Range.Text = "52-44="

$t.Cell(4, 1).Range.Text = "1+78="
$t.Cell(4, 2).Range.Text = "35+13="
$t.Cell(4, 3).Range.Text = "19+11="
$t.Cell(4, 4).Range.Text = "44-33="
$t.Cell(4, 5).Range.Text = "90-70="

$t.Cell(5, 1).Range.Text = "73-37="
$t.Cell(5, 2).Range.Text = "67-53="
$t.Cell(5, 3).Range.Text = "72-62="
$t.Cell(5, 4).Range.Text = "61-41="
$t.Cell(5, 5).Range.Text = "67+8="

$t.Cell(6, 1).Range.Text = "80-10="
$t.Cell(6, 2).Range.Text = "39-31="
$t.Cell(6, 3).Range.Text = "68+19="
$t.Cell(6, 4).Range.Text = "9+0="
$t.Cell(6, 5).Range.Text = "53-22="

$t.Cell(7, 1).Range.Text = "65-7="
$t.Cell(7, 2).Range.Text = "55-13="
$t.Cell(7, 3).Range.Text = "8+50="
$t.Cell(7, 4).Range.Text = "31+68="
$t.Cell(7, 5).Range.Text = "47-21="

$t.Cell(8, 1).Range.Text = "28+16="
$t.Cell(8, 2).Range.Text = "99-33="
$t.Cell(8, 3).Range.Text = "86-24="
$t.Cell(8, 4).Range.Text = "62-23="
$t.Cell(8, 5).Range.Text = "37-2="

$t.Cell(9, 1).Range.Text = "78-69="
$t.Cell(9, 2).Range.Text = "66-66="
$t.Cell(9, 3).Range.Text = "23-5="
$t.Cell(9, 4).Range.Text = "54+7="
$t.Cell(9, 5).Range.Text = "30+62="

$t.Cell(10, 1).Range.Text = "0+66="
$t.Cell(10, 2).Range.Text = "3+20="
$t.Cell(10, 3).Range.Text = "46-15="
$t.Cell(10, 4).Range.Text = "64+1="
$t.Cell(10, 5).Range.Text = "95-31="

$t.Cell(11, 1).Range.Text = "29-0="
$t.Cell(11, 2).Range.Text = "72-49="
$t.Cell(11, 3).Range.Text = "5+56="
$t.Cell(11, 4).Range.Text = "48+34="
$t.Cell(11, 5).Range.Text = "1+29="

$t.Cell(12, 1).Range.Text = "66-45="
$t.Cell(12, 2).Range.Text = "97-19="
$t.Cell(12, 3).Range.Text = "15+63="
$t.Cell(12, 4).Range.Text = "74+12="
$t.Cell(12, 5).Range.Text = "9+4="

$t.Cell(13, 1).Range.Text = "54-2="
$t.Cell(13, 2).Range.Text = "81-1="
$t.Cell(13, 3).Range.Text = "45+3="
$t.Cell(13, 4).Range.Text = "6+40="
$t.Cell(13, 5).Range.Text = "96-80="

$t.Cell(14, 1).Range.Text = "33-16="
$t.Cell(14, 2).Range.Text = "68+12="
$t.Cell(14, 3).Range.Text = "6+88="
$t.Cell(14, 4).Range.Text = "64-5="
$t.Cell(14, 5).Range.Text = "8+84="

$t.Cell(15, 1).Range.Text = "55-29="
$t.Cell(15, 2).Range.Text = "13+64="
$t.Cell(15, 3).Range.Text = "26+42="
$t.Cell(15, 4).Range.Text = "73-9="
$t.Cell(15, 5).Range.Text = "7+92="

$t.Cell(16, 1).Range.Text = "6+31="
$t.Cell(16, 2).Range.Text = "67+30="
$t.Cell(16, 3).Range.Text = "81-61="
$t.Cell(16, 4).Range.Text = "37+36="
$t.Cell(16, 5).Range.Text = "56-53="

$t.Cell(17, 1).Range.Text = "17+15="
$t.Cell(17, 2).Range.Text = "92-51="
$t.Cell(17, 3).Range.Text = "76+1="
$t.Cell(17, 4).Range.Text = "70-30="
$t.Cell(17, 5).Range.Text = "64+17="

$t.Cell(18, 1).Range.Text = "92+6="
$t.Cell(18, 2).Range.Text = "84-13="
$t.Cell(18, 3).Range.Text = "6+6="
$t.Cell(18, 4).Range.Text = "30+4="
$t.Cell(18, 5).Range.Text = "4+8="

$t.Cell(19, 1).Range.Text = "85-2="
$t.Cell(19, 2).Range.Text = "86-9="
$t.Cell(19, 3).Range.Text = "48-4="
$t.Cell(19, 4).Range.Text = "28-2="
$t.Cell(19, 5).Range.Text = "61+18="

$t.Cell(20, 1).Range.Text = "36+26="
$t.Cell(20, 2).Range.Text = "47-40="
$t.Cell(20, 3).Range.Text = "10+59="
$t.Cell(20, 4).Range.Text = "11+70="
$t.Cell(20, 5).Range.Text = "13+46="

